$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text columns (Coin name, Link) ---
$textUpdates = @(
    @{Cell="B36"; Value="ARBITRUM"},
    @{Cell="B37"; Value="Celestia"},
    @{Cell="B44"; Value="Monero"},
    @{Cell="B45"; Value="EnergySwap"},
    @{Cell="C36"; Value="https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"},
    @{Cell="C37"; Value="https://coinranking.com/coin/YQcD0lBl7+celestia-tia"},
    @{Cell="C44"; Value="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"},
    @{Cell="C45"; Value="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"}
)
foreach ($item in $textUpdates) {
    $ws.Range($item.Cell).Value = $item.Value
}

# --- Volume/percent column (E) - already text-safe due to leading/trailing spaces ---
$eUpdates = @(
    @{Cell="E2"; Value="  +4.99%  "},
    @{Cell="E3"; Value="  +5.22%  "},
    @{Cell="E4"; Value="  +0.00%  "},
    @{Cell="E5"; Value="  +4.67%  "},
    @{Cell="E6"; Value="  +2.60%  "},
    @{Cell="E7"; Value="  +4.69%  "},
    @{Cell="E8"; Value="  -0.02%  "},
    @{Cell="E9"; Value="  +5.00%  "},
    @{Cell="E10"; Value="  +5.50%  "},
    @{Cell="E11"; Value="  +5.22%  "},
    @{Cell="E12"; Value="  +0.49%  "},
    @{Cell="E14"; Value="  -0.16%  "},
    @{Cell="E15"; Value="  +5.10%  "},
    @{Cell="E16"; Value="  +5.43%  "},
    @{Cell="E17"; Value="  +4.69%  "},
    @{Cell="E18"; Value="  +2.04%  "},
    @{Cell="E19"; Value="  +9.84%  "},
    @{Cell="E20"; Value="  +4.93%  "},
    @{Cell="E21"; Value="  -1.16%  "},
    @{Cell="E23"; Value="  +2.87%  "},
    @{Cell="E24"; Value="  +1.44%  "},
    @{Cell="E25"; Value="  +7.32%  "},
    @{Cell="E26"; Value="  +2.26%  "},
    @{Cell="E27"; Value="  -0.02%  "},
    @{Cell="E28"; Value="  +0.12%  "},
    @{Cell="E29"; Value="  +0.86%  "},
    @{Cell="E30"; Value="  +1.81%  "},
    @{Cell="E31"; Value="  +0.08%  "},
    @{Cell="E32"; Value="  +1.43%  "},
    @{Cell="E33"; Value="  +4.25%  "},
    @{Cell="E34"; Value="  -0.41%  "},
    @{Cell="E35"; Value="  -0.04%  "},
    @{Cell="E36"; Value="  +3.57%  "},
    @{Cell="E37"; Value="  -1.08%  "},
    @{Cell="E38"; Value="  +0.41%  "},
    @{Cell="E39"; Value="  +3.45%  "},
    @{Cell="E40"; Value="  +9.97%  "},
    @{Cell="E41"; Value="  +26.92%  "},
    @{Cell="E42"; Value="  +3.64%  "},
    @{Cell="E43"; Value="  +3.16%  "},
    @{Cell="E44"; Value="  -3.07%  "},
    @{Cell="E45"; Value="  -2.82%  "},
    @{Cell="E47"; Value="  +0.69%  "},
    @{Cell="E48"; Value="  +0.48%  "},
    @{Cell="E50"; Value="  -0.96%  "},
    @{Cell="E51"; Value="  +14.15%  "}
)
foreach ($item in $eUpdates) {
    $ws.Range($item.Cell).Value = $item.Value
}

# --- Price column (D) - force text to preserve exact formatting (leading/trailing zeros, multi-dot numbers) ---
$dUpdates = @(
    @{Cell="D2"; Value="52.048.85"},
    @{Cell="D3"; Value="2.780.45"},
    @{Cell="D5"; Value="342.22"},
    @{Cell="D6"; Value="115.43"},
    @{Cell="D9"; Value="0.577"},
    @{Cell="D10"; Value="41.84"},
    @{Cell="D12"; Value="20.07"},
    @{Cell="D14"; Value="7.62"},
    @{Cell="D15"; Value="3.216.52"},
    @{Cell="D16"; Value="2.778.46"},
    @{Cell="D17"; Value="51.875.45"},
    @{Cell="D18"; Value="0.876"},
    @{Cell="D20"; Value="7.01"},
    @{Cell="D21"; Value="13.23"},
    @{Cell="D22"; Value="0.0₃0977"},
    @{Cell="D23"; Value="276.25"},
    @{Cell="D24"; Value="69.95"},
    @{Cell="D25"; Value="2.76"},
    @{Cell="D26"; Value="26.68"},
    @{Cell="D27"; Value="0.999"},
    @{Cell="D30"; Value="0.141"},
    @{Cell="D31"; Value="34.67"},
    @{Cell="D32"; Value="50.31"},
    @{Cell="D36"; Value="2.10"},
    @{Cell="D37"; Value="19.00"},
    @{Cell="D39"; Value="3.22"},
    @{Cell="D40"; Value="0.0383"},
    @{Cell="D44"; Value="125.84"},
    @{Cell="D45"; Value="23.06"},
    @{Cell="D46"; Value="2.068.17"},
    @{Cell="D51"; Value="0.880"}
)
foreach ($item in $dUpdates) {
    $r = $ws.Range($item.Cell)
    $r.NumberFormat = "@"
    $r.Value = $item.Value
    $r.ClearFormats()
}
